# Auto-generated: applies scheduled-runner price/profit updates to the
# Shiva_Profits workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each block updates the currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) for specific leve rows based on refreshed market data.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 75002500
$ws.Range("I32").Value = 200000000
$ws.Range("J32").Value = 33336662
$ws.Range("K32").Value = 200000000
$ws.Range("L32").Value = 33336662
$ws.Range("M32").Value = -199999674
$ws.Range("N32").Value = -33337314
$ws.Range("H69").Value = 14220.385
$ws.Range("I69").Value = 11186.6
$ws.Range("K69").Value = 33559.8
$ws.Range("M69").Value = -32685.8
$ws.Range("H70").Value = 1242.1666
$ws.Range("I70").Value = 2102.2856
$ws.Range("J70").Value = 888
$ws.Range("K70").Value = 6306.8568
$ws.Range("L70").Value = 2664
$ws.Range("M70").Value = -6036.8568
$ws.Range("N70").Value = -3204
$ws.Range("H72").Value = 14220.385
$ws.Range("I72").Value = 11186.6
$ws.Range("K72").Value = 100679.4
$ws.Range("M72").Value = -96311.40000000001
$ws.Range("H73").Value = 1242.1666
$ws.Range("I73").Value = 2102.2856
$ws.Range("J73").Value = 888
$ws.Range("K73").Value = 6306.8568
$ws.Range("L73").Value = 2664
$ws.Range("M73").Value = -5370.8568
$ws.Range("N73").Value = -4536
$ws.Range("H92").Value = 1199.7727
$ws.Range("I92").Value = 1115.9286
$ws.Range("J92").Value = 1346.5
$ws.Range("K92").Value = 1115.9286
$ws.Range("L92").Value = 1346.5
$ws.Range("M92").Value = 132.0714
$ws.Range("N92").Value = -3842.5
$ws.Range("H97").Value = 2182.4167
$ws.Range("J97").Value = 2182.4167
$ws.Range("L97").Value = 6547.250100000001
$ws.Range("N97").Value = -7539.250100000001
$ws.Range("H100").Value = 1247.7333
$ws.Range("I100").Value = 1105.2727
$ws.Range("K100").Value = 1105.2727
$ws.Range("M100").Value = -564.2727
$ws.Range("H137").Value = 4175.1206
$ws.Range("I137").Value = 4302.213
$ws.Range("J137").Value = 3632.0908
$ws.Range("K137").Value = 12906.639
$ws.Range("L137").Value = 10896.2724
$ws.Range("M137").Value = -10356.639
$ws.Range("N137").Value = -15996.2724

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5887.557
$ws.Range("I32").Value = 4003.2466
$ws.Range("K32").Value = 4003.2466
$ws.Range("M32").Value = -3716.2466
$ws.Range("H61").Value = 5619
$ws.Range("I61").Value = 5481.0527
$ws.Range("K61").Value = 5481.0527
$ws.Range("M61").Value = -5269.0527
$ws.Range("H122").Value = 17649.592
$ws.Range("I122").Value = 2839.7778
$ws.Range("K122").Value = 8519.3334
$ws.Range("M122").Value = -6069.3334
$ws.Range("H136").Value = 5619
$ws.Range("I136").Value = 5481.0527
$ws.Range("K136").Value = 16443.1581
$ws.Range("M136").Value = -13893.1581

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 150047.81
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 150047.81
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 150047.81
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -160327.81

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 200
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 200
$ws.Range("N4").Value = -424
$ws.Range("H31").Value = 1744.4147
$ws.Range("I31").Value = 1573.6578
$ws.Range("J31").Value = 3907.3333
$ws.Range("K31").Value = 1573.6578
$ws.Range("L31").Value = 3907.3333
$ws.Range("M31").Value = -1278.6578
$ws.Range("N31").Value = -4497.3333
$ws.Range("H34").Value = 1744.4147
$ws.Range("I34").Value = 1573.6578
$ws.Range("J34").Value = 3907.3333
$ws.Range("K34").Value = 1573.6578
$ws.Range("L34").Value = 3907.3333
$ws.Range("M34").Value = -1371.6578
$ws.Range("N34").Value = -4311.3333
$ws.Range("H51").Value = 57500
$ws.Range("J51").Value = 57500
$ws.Range("L51").Value = 57500
$ws.Range("N51").Value = -58972
$ws.Range("H58").Value = 1994.2307
$ws.Range("I58").Value = 990.625
$ws.Range("J58").Value = 3600
$ws.Range("K58").Value = 990.625
$ws.Range("L58").Value = 3600
$ws.Range("M58").Value = -787.625
$ws.Range("N58").Value = -4006
$ws.Range("H61").Value = 57500
$ws.Range("J61").Value = 57500
$ws.Range("L61").Value = 57500
$ws.Range("N61").Value = -58196
$ws.Range("H68").Value = 45890.184
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 45890.184
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H107").Value = 1790.75
$ws.Range("I107").Value = 1766
$ws.Range("K107").Value = 1766
$ws.Range("M107").Value = 154
$ws.Range("H120").Value = 26599.334
$ws.Range("J120").Value = 26599.334
$ws.Range("L120").Value = 26599.334
$ws.Range("N120").Value = -33857.334
$ws.Range("H121").Value = 59999.5
$ws.Range("J121").Value = 59999.5
$ws.Range("L121").Value = 59999.5
$ws.Range("N121").Value = -62619.5
$ws.Range("H122").Value = 6618.963
$ws.Range("I122").Value = 11268
$ws.Range("K122").Value = 33804
$ws.Range("M122").Value = -31354
$ws.Range("H133").Value = 48612.332
$ws.Range("J133").Value = 48612.332
$ws.Range("L133").Value = 48612.332
$ws.Range("N133").Value = -53672.332
$ws.Range("H136").Value = 1994.2307
$ws.Range("I136").Value = 990.625
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 2971.875
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = -421.875
$ws.Range("N136").Value = -15900
$ws.Range("H141").Value = 131709.7
$ws.Range("J141").Value = 140229.3
$ws.Range("L141").Value = 140229.3
$ws.Range("N141").Value = -150589.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 30561.848
$ws.Range("I4").Value = 131.14815
$ws.Range("K4").Value = 393.44445
$ws.Range("M4").Value = -281.44445
$ws.Range("H9").Value = 91820850
$ws.Range("I9").Value = 5000500.5
$ws.Range("J9").Value = 111114260
$ws.Range("K9").Value = 15001501.5
$ws.Range("L9").Value = 333342780
$ws.Range("M9").Value = -15001277.5
$ws.Range("N9").Value = -333343228

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 222.05263
$ws.Range("J97").Value = 259.8889
$ws.Range("L97").Value = 259.8889
$ws.Range("N97").Value = -1251.8889
$ws.Range("H102").Value = 9026.526
$ws.Range("I102").Value = 9880.333000000001
$ws.Range("J102").Value = 7562.857
$ws.Range("K102").Value = 9880.333000000001
$ws.Range("L102").Value = 7562.857
$ws.Range("M102").Value = -8258.333000000001
$ws.Range("N102").Value = -10806.857
$ws.Range("H107").Value = 1150.05
$ws.Range("J107").Value = 452.14285
$ws.Range("L107").Value = 452.14285
$ws.Range("N107").Value = -4292.14285
$ws.Range("H109").Value = 32283.5
$ws.Range("J109").Value = 32283.5
$ws.Range("L109").Value = 32283.5
$ws.Range("N109").Value = -34363.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7996.1665
$ws.Range("I40").Value = 7997.25
$ws.Range("J40").Value = 7994
$ws.Range("K40").Value = 7997.25
$ws.Range("L40").Value = 7994
$ws.Range("M40").Value = -7861.25
$ws.Range("N40").Value = -8266
$ws.Range("H122").Value = 3761.762
$ws.Range("I122").Value = 2818.375
$ws.Range("K122").Value = 8455.125
$ws.Range("M122").Value = -6005.125
$ws.Range("H139").Value = 134635.67
$ws.Range("J139").Value = 127453.5
$ws.Range("L139").Value = 127453.5
$ws.Range("N139").Value = -137733.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 29421220
$ws.Range("I62").Value = 50011884
$ws.Range("J62").Value = 5987
$ws.Range("K62").Value = 50011884
$ws.Range("L62").Value = 5987
$ws.Range("M62").Value = -50011260
$ws.Range("N62").Value = -7235
$ws.Range("H65").Value = 29421220
$ws.Range("I65").Value = 50011884
$ws.Range("J65").Value = 5987
$ws.Range("K65").Value = 250059420
$ws.Range("L65").Value = 29935
$ws.Range("M65").Value = -250056300
$ws.Range("N65").Value = -36175
$ws.Range("H96").Value = 2887.2354
$ws.Range("I96").Value = 2980.5
$ws.Range("J96").Value = 2754
$ws.Range("K96").Value = 2980.5
$ws.Range("L96").Value = 2754
$ws.Range("M96").Value = -1607.5
$ws.Range("N96").Value = -5500
$ws.Range("H113").Value = 745.375
$ws.Range("J113").Value = 473.7
$ws.Range("L113").Value = 1421.1
$ws.Range("N113").Value = -5761.1
$ws.Range("H118").Value = 107499.75
$ws.Range("J118").Value = 107499.75
$ws.Range("L118").Value = 107499.75
$ws.Range("N118").Value = -110813.75
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H136").Value = 1661.7451
$ws.Range("I136").Value = 1686.1957
$ws.Range("J136").Value = 1436.8
$ws.Range("K136").Value = 5058.5871
$ws.Range("L136").Value = 4310.4
$ws.Range("M136").Value = -2508.5871
$ws.Range("N136").Value = -9410.4
$ws.Range("H139").Value = 69899.14
$ws.Range("J139").Value = 69899.14
$ws.Range("L139").Value = 69899.14
$ws.Range("N139").Value = -80179.14

